# Updates the cryptos list (Price / Volume(1h) columns, plus a few
# Coin/Link cells where rows shifted) to match the refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (prevents Excel from reinterpreting
# numeric-looking strings like "1.00" or "51.544.64" as numbers), and then
# clear the style back to Normal so no stray number-format/quote-prefix
# styling gets attached to the cell.
function Set-TextCell($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "51.544.64"
$ws.Range("E2").Value = "  -1.28%  "

Set-TextCell $ws.Range("D3") "2.919.24"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.14%  "

Set-TextCell $ws.Range("D5") "350.54"
$ws.Range("E5").Value = "  -0.91%  "

Set-TextCell $ws.Range("D6") "106.18"
$ws.Range("E6").Value = "  -6.86%  "

Set-TextCell $ws.Range("D7") "0.554"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("E8").Value = "  -0.02%  "

Set-TextCell $ws.Range("D9") "0.608"
$ws.Range("E9").Value = "  -2.47%  "

Set-TextCell $ws.Range("D10") "37.50"
$ws.Range("E10").Value = "  -5.62%  "

$ws.Range("E11").Value = "  +0.94%  "

Set-TextCell $ws.Range("D12") "0.0850"
$ws.Range("E12").Value = "  -2.16%  "

Set-TextCell $ws.Range("D13") "18.94"
$ws.Range("E13").Value = "  -4.44%  "

Set-TextCell $ws.Range("D14") "3.383.30"
$ws.Range("E14").Value = "  +0.42%  "

Set-TextCell $ws.Range("D15") "7.64"
$ws.Range("E15").Value = "  -1.39%  "

Set-TextCell $ws.Range("D16") "2.909.49"
$ws.Range("E16").Value = "  +0.03%  "

Set-TextCell $ws.Range("D17") "0.960"
$ws.Range("E17").Value = "  -2.73%  "

Set-TextCell $ws.Range("D18") "51.550.75"
$ws.Range("E18").Value = "  -1.41%  "

Set-TextCell $ws.Range("D19") "3.40"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("E20").Value = "  -3.45%  "

Set-TextCell $ws.Range("D21") "13.34"
$ws.Range("E21").Value = "  -5.40%  "

Set-TextCell $ws.Range("D22") "0.0₃0956"

Set-TextCell $ws.Range("D23") "68.82"
$ws.Range("E23").Value = "  -3.17%  "

Set-TextCell $ws.Range("D24") "260.85"
$ws.Range("E24").Value = "  -3.15%  "

Set-TextCell $ws.Range("D25") "2.71"
$ws.Range("E25").Value = "  -3.70%  "

Set-TextCell $ws.Range("D26") "0.171"
$ws.Range("E26").Value = "  -5.74%  "

Set-TextCell $ws.Range("D27") "26.42"
$ws.Range("E27").Value = "  -1.34%  "

Set-TextCell $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.06%  "

Set-TextCell $ws.Range("D29") "7.41"
$ws.Range("E29").Value = "  +8.91%  "

$ws.Range("E30").Value = "  -0.43%  "

Set-TextCell $ws.Range("D31") "10.18"
$ws.Range("E31").Value = "  -4.49%  "

$ws.Range("E32").Value = "  -5.29%  "

Set-TextCell $ws.Range("D33") "35.58"
$ws.Range("E33").Value = "  -5.27%  "

Set-TextCell $ws.Range("D34") "5.90"
$ws.Range("E34").Value = "  -3.31%  "

Set-TextCell $ws.Range("D35") "51.01"
$ws.Range("E35").Value = "  -3.88%  "

Set-TextCell $ws.Range("D36") "0.998"
$ws.Range("E36").Value = "  -0.18%  "

Set-TextCell $ws.Range("D37") "0.0423"
$ws.Range("E37").Value = "  -6.26%  "

Set-TextCell $ws.Range("B38") "LidoDAOToken"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws.Range("D38") "3.12"
$ws.Range("E38").Value = "  -5.91%  "

Set-TextCell $ws.Range("B39") "ARBITRUM"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D39") "1.95"
$ws.Range("E39").Value = "  -4.41%  "

Set-TextCell $ws.Range("B40") "Celestia"
Set-TextCell $ws.Range("C40") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell $ws.Range("D40") "17.60"
$ws.Range("E40").Value = "  -6.53%  "

Set-TextCell $ws.Range("B41") "Stacks"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D41") "2.63"
$ws.Range("E41").Value = "  -3.98%  "

Set-TextCell $ws.Range("B42") "Stellar"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D42") "0.116"
$ws.Range("E42").Value = "  -1.40%  "

Set-TextCell $ws.Range("B43") "EnergySwap"
Set-TextCell $ws.Range("C43") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D43") "22.60"
$ws.Range("E43").Value = "  -2.05%  "

Set-TextCell $ws.Range("D44") "119.24"
$ws.Range("E44").Value = "  +1.17%  "

Set-TextCell $ws.Range("D45") "2.15"
$ws.Range("E45").Value = "  -0.85%  "

Set-TextCell $ws.Range("D46") "2.45"
$ws.Range("E46").Value = "  -3.46%  "

Set-TextCell $ws.Range("D47") "2.091.78"
$ws.Range("E47").Value = "  -4.03%  "

$ws.Range("E48").Value = "  -6.53%  "

Set-TextCell $ws.Range("B49") "RocketPoolETH"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws.Range("D49") "3.213.05"
$ws.Range("E49").Value = "  +0.45%  "

Set-TextCell $ws.Range("B50") "TheGraph"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell $ws.Range("D50") "0.236"
$ws.Range("E50").Value = "  -9.43%  "

Set-TextCell $ws.Range("B51") "BEAM"
Set-TextCell $ws.Range("C51") "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextCell $ws.Range("D51") "0.0337"
$ws.Range("E51").Value = "  -5.02%  "

